$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 32-33: dataset links (Hyperlink style)
$ws.Range("A32").Value = "https://data.census.gov/table"
$ws.Hyperlinks.Add($ws.Range("A32"), "https://data.census.gov/table")
$ws.Range("A32").Style = "Hyperlink"

$ws.Range("A33").Value = "https://www.kaggle.com/datasets/utkarshxy/who-worldhealth-statistics-2020-complete"
$ws.Hyperlinks.Add($ws.Range("A33"), "https://www.kaggle.com/datasets/utkarshxy/who-worldhealth-statistics-2020-complete")
$ws.Range("A33").Style = "Hyperlink"

# Row 34-35: empty cells carrying the Hyperlink style (spacer rows)
$ws.Range("A34").Style = "Hyperlink"
$ws.Range("A35").Style = "Hyperlink"

# Row 37: plain-text section header (no hyperlink style)
$ws.Range("A37").Value = "Example capstones, reports, articles, etc."

# Rows 38-41: capstone/article example links
$ws.Range("A38").Value = "https://jscholarship.library.jhu.edu/bitstream/handle/1774.2/61821/Wallace,%20Robyn.pdf?sequence=1"
$ws.Hyperlinks.Add($ws.Range("A38"), "https://jscholarship.library.jhu.edu/bitstream/handle/1774.2/61821/Wallace,%20Robyn.pdf?sequence=1")
$ws.Range("A38").Style = "Hyperlink"

$ws.Range("A39").Value = "https://pubmed.ncbi.nlm.nih.gov/36276352/"
$ws.Hyperlinks.Add($ws.Range("A39"), "https://pubmed.ncbi.nlm.nih.gov/36276352/")
$ws.Range("A39").Style = "Hyperlink"

$ws.Range("A40").Value = "https://digitalcommons.csumb.edu/cgi/viewcontent.cgi?article=1849&context=caps_thes_all"
$ws.Hyperlinks.Add($ws.Range("A40"), "https://digitalcommons.csumb.edu/cgi/viewcontent.cgi?article=1849&context=caps_thes_all")
$ws.Range("A40").Style = "Hyperlink"

$ws.Range("A41").Value = "https://scholarworks.gsu.edu/cgi/viewcontent.cgi?article=1022&context=iph_capstone"
$ws.Hyperlinks.Add($ws.Range("A41"), "https://scholarworks.gsu.edu/cgi/viewcontent.cgi?article=1022&context=iph_capstone")
$ws.Range("A41").Style = "Hyperlink"

# Rows 42-44: race/health policy data links
$ws.Range("A42").Value = "https://www.brookings.edu/blog/usc-brookings-schaeffer-on-health-policy/2020/02/19/there-are-clear-race-based-inequalities-in-health-insurance-and-health-outcomes/"
$ws.Hyperlinks.Add($ws.Range("A42"), "https://www.brookings.edu/blog/usc-brookings-schaeffer-on-health-policy/2020/02/19/there-are-clear-race-based-inequalities-in-health-insurance-and-health-outcomes/")
$ws.Range("A42").Style = "Hyperlink"

$ws.Range("A43").Value = "https://www.kff.org/racial-equity-and-health-policy/report/key-data-on-health-and-health-care-by-race-and-ethnicity/"
$ws.Hyperlinks.Add($ws.Range("A43"), "https://www.kff.org/racial-equity-and-health-policy/report/key-data-on-health-and-health-care-by-race-and-ethnicity/")
$ws.Range("A43").Style = "Hyperlink"

$ws.Range("A44").Value = "https://www.kff.org/report-section/key-data-on-health-and-health-care-by-race-ethnicity-methodology/"
$ws.Hyperlinks.Add($ws.Range("A44"), "https://www.kff.org/report-section/key-data-on-health-and-health-care-by-race-ethnicity-methodology/")
$ws.Range("A44").Style = "Hyperlink"

# Row 45: kaggle happiness notebook link
$ws.Range("A45").Value = "https://www.kaggle.com/code/noobiedatascientist/explaining-happiness"
$ws.Hyperlinks.Add($ws.Range("A45"), "https://www.kaggle.com/code/noobiedatascientist/explaining-happiness")
$ws.Range("A45").Style = "Hyperlink"

# Final selection, matching the saved view state in the workbook
$ws.Range("A35").Select() | Out-Null
